$wb = $excel.ActiveWorkbook

# --- Step 1: select B1 on the Content sheet (matches diff selection change) ---
$wsContent = $wb.Worksheets.Item(2)
$wsContent.Range("B1").Select()

# --- Step 2: add the new "Messages" worksheet after "Content" ---
$wsMessages = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$wsMessages.Name = "Messages"

# --- Step 3: write cell values in the exact order required to reproduce the
#     shared-string table order 58..106 from the target workbook ---
$wsMessages.Range("A2").Value = "en"
$wsMessages.Range("A3").Value = "nl"
$wsMessages.Range("A5").Value = "es"
$wsMessages.Range("A4").Value = "aw"
$wsMessages.Range("B1").Value = "toolbarTitle"
$wsMessages.Range("C1").Value = "listTitle"
$wsMessages.Range("D1").Value = "loadingText"
$wsMessages.Range("E1").Value = "submitButton"
$wsMessages.Range("C2").Value = "What did you throw away?"
$wsMessages.Range("B2").Value = "Let's collect data!"
$wsMessages.Range("D2").Value = "Loading application..."
$wsMessages.Range("E2").Value = "Submit"
$wsMessages.Range("B3").Value = "Laten we data verzamelen!"
$wsMessages.Range("C3").Value = "Wat heeft u weggegooid?"
$wsMessages.Range("E3").Value = "Verstuur"
$wsMessages.Range("F1").Value = "defaultThankYouMessage"
$wsMessages.Range("G1").Value = "co2MessagePart1"
$wsMessages.Range("H1").Value = "co2MessagePart2"
$wsMessages.Range("I1").Value = "askFeedback"
$wsMessages.Range("J1").Value = "clickHere"
$wsMessages.Range("K1").Value = "normalBin"
$wsMessages.Range("F2").Value = "Thank you for registering your waste!"
$wsMessages.Range("G2").Value = "This recycle bin already saved"
$wsMessages.Range("H2").Value = "kg of CO2!"
$wsMessages.Range("I2").Value = "Would you be so kind to give some feedback on the app?"
$wsMessages.Range("J2").Value = "Click here"
$wsMessages.Range("K2").Value = "Normal bin"
$wsMessages.Range("D3").Value = "Applicatie aan het laden..."
$wsMessages.Range("F3").Value = "Bedankt om uw afval te registreren!"
$wsMessages.Range("G3").Value = "Deze vuilbak heeft al"
$wsMessages.Range("H3").Value = "kg CO2 bespaard!"
$wsMessages.Range("I3").Value = "Zou u aub wat feedback willen geven over de app?"
$wsMessages.Range("J3").Value = "Klik hier"
$wsMessages.Range("K3").Value = "Normale vuilbak"
$wsMessages.Range("A1").Value = "language"
$wsMessages.Range("B4").Value = "Papiamento!"
$wsMessages.Range("B5").Value = "Español!"
$wsMessages.Range("L1").Value = "recyclingBin"
$wsMessages.Range("M1").Value = "nonRecyclingBin"
$wsMessages.Range("N1").Value = "hasToBeRecyclingBin"
$wsMessages.Range("O1").Value = "hasToBeNonRecyclingBin"
$wsMessages.Range("L2").Value = "Recycling bin"
$wsMessages.Range("M2").Value = "Non-recycling bin"
$wsMessages.Range("N2").Value = "This item belongs in the recycling bin"
$wsMessages.Range("O2").Value = "This item belongs in the non-recycling bin"
$wsMessages.Range("L3").Value = "Recyclingvuilbak"
$wsMessages.Range("M3").Value = "Niet-recyclingvuilbak"
$wsMessages.Range("N3").Value = "Dit item behoort in de recyclingvuilbak"
$wsMessages.Range("O3").Value = "Dit item behoort in de niet-recyclingvuilbak"

# --- Step 4: fill in the remaining (duplicate-valued) cells ---
$wsMessages.Range("C4").Value = "What did you throw away?"
$wsMessages.Range("D4").Value = "Loading application..."
$wsMessages.Range("E4").Value = "Submit"
$wsMessages.Range("F4").Value = "Thank you for registering your waste!"
$wsMessages.Range("G4").Value = "This recycle bin already saved"
$wsMessages.Range("H4").Value = "kg of CO2!"
$wsMessages.Range("I4").Value = "Would you be so kind to give some feedback on the app?"
$wsMessages.Range("J4").Value = "Click here"
$wsMessages.Range("K4").Value = "Normal bin"
$wsMessages.Range("L4").Value = "Recyclingvuilbak"
$wsMessages.Range("M4").Value = "Niet-recyclingvuilbak"
$wsMessages.Range("N4").Value = "Dit item behoort in de recyclingvuilbak"
$wsMessages.Range("O4").Value = "Dit item behoort in de niet-recyclingvuilbak"
$wsMessages.Range("C5").Value = "What did you throw away?"
$wsMessages.Range("D5").Value = "Loading application..."
$wsMessages.Range("E5").Value = "Submit"
$wsMessages.Range("F5").Value = "Thank you for registering your waste!"
$wsMessages.Range("G5").Value = "This recycle bin already saved"
$wsMessages.Range("H5").Value = "kg of CO2!"
$wsMessages.Range("I5").Value = "Would you be so kind to give some feedback on the app?"
$wsMessages.Range("J5").Value = "Click here"
$wsMessages.Range("K5").Value = "Normal bin"
$wsMessages.Range("L5").Value = "Recyclingvuilbak"
$wsMessages.Range("M5").Value = "Niet-recyclingvuilbak"
$wsMessages.Range("N5").Value = "Dit item behoort in de recyclingvuilbak"
$wsMessages.Range("O5").Value = "Dit item behoort in de niet-recyclingvuilbak"

# --- Step 5: column widths (best achievable approximation; engine rounds to 1/6 char) ---
$wsMessages.Range("A1:B1").EntireColumn.ColumnWidth = 18.333333333333332
$wsMessages.Range("C1:D1").EntireColumn.ColumnWidth = 29.666666666666668
$wsMessages.Range("E1:E1").EntireColumn.ColumnWidth = 19.166666666666668
$wsMessages.Range("F1:K1").EntireColumn.ColumnWidth = 36.333333333333336
$wsMessages.Range("L1:L1").EntireColumn.ColumnWidth = 24.666666666666668
$wsMessages.Range("M1:M1").EntireColumn.ColumnWidth = 22.333333333333332
$wsMessages.Range("N1:N1").EntireColumn.ColumnWidth = 34.0
$wsMessages.Range("O1:O1").EntireColumn.ColumnWidth = 37.166666666666664
$wsMessages.Range("P1:P1").EntireColumn.ColumnWidth = 22.333333333333332

# --- Step 6: selection on Messages sheet ---
$wsMessages.Range("K6").Select()

Write-Output "done"
